$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns we touch so Excel does not
# reinterpret dotted/percentage strings as numbers (matching original inlineStr text cells).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.500.05"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.877.85"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").Value = "1.023"
$ws.Range("E4").Value = "  +1.64%  "
$ws.Range("D5").Value = "317.28"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "1.021"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("D7").Value = "0.5148"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "0.3946"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").Value = "0.08338"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "1.118"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "42.07"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").Value = "6.248"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.866.64"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "20.46"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "7.242"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "1.023"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "91.44"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "0.06759"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "17.72"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.021"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").Value = "5.972"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").Value = "28.529.36"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").Value = "2.268"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "2.075.23"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "161.71"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("D28").Value = "20.78"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "2.381"
$ws.Range("E29").Value = "  -4.66%  "
$ws.Range("D30").Value = "127.32"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").Value = "0.1053"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").Value = "1.037"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "5.842"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "3.659"
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("D35").Value = "0.02438"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "0.06505"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "9.144"
$ws.Range("E37").Value = "  -5.63%  "
$ws.Range("D38").Value = "0.2188"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").Value = "1.254"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").Value = "1.191"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").Value = "0.6446"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").Value = "5.006"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "11.20"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").Value = "0.6037"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").Value = "13.02"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "3.715"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").Value = "1.258"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "1.215"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("D50").Value = "122.14"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "0.06878"
$ws.Range("E51").Value = "  -0.23%  "
